# Customers.xlsx — "Added functionalities to all admin use classes"
#
# Sheet1 ("Customers"): refresh the two sample-customer rows with new
# banking data, move the sample e-mail addresses from the abc.com domain
# to gmail.com, and drop the hyperlink that used to sit on the first
# sample row (only the second row keeps a live mailto: link now).
#
# Sheet2 ("Admin"): add a second sample company ("DEF") below the
# existing "ABC" admin row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet1: update the two customer rows
# ---------------------------------------------------------------------

$ws1.Range("C3").Value = "sample1@gmail.com"
$ws1.Range("D3").Value = 1615
$ws1.Range("E3").Value = 100136144
$ws1.Range("F3").Value = 1500

$ws1.Range("C4").Value = "sample2@gmail.com"
$ws1.Range("D4").Value = 9766
$ws1.Range("E4").Value = 100195188
$ws1.Range("F4").Value = 1500

# Drop both existing hyperlinks, then re-add a single live link on the
# second row only (first row keeps the Hyperlink look but is no longer
# a clickable link).
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("C4"), "mailto:sample2@gmail.com")

# Widen column E (AccNo) now that the account numbers are much longer.
$ws1.Columns.Item(5).ColumnWidth = 9.166666666666666

# ---------------------------------------------------------------------
# Sheet2: add the new "DEF" admin/company row
# ---------------------------------------------------------------------

$ws2.Range("A4").Value = "DEF company"
$ws2.Range("B4").Value = "DEF"
$ws2.Range("C4").Value = "def@gmail.com"
$ws2.Range("D4").Value = 2348
$ws2.Range("E4").Value = 900141676
$ws2.Range("F4").Value = 0

$ws2.Range("E3").Value = 900113678

$ws2.Columns.Item(5).ColumnWidth = 10.166666666666666

$ws2.Range("C9").Select() | Out-Null
